$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header): add P1=14, Q1=15, matching the bold/bordered style of B1:O1 ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the formatting from O1 (last existing styled header cell) onto the new P1:Q1 cells
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows 2-25: update I, K, M, O columns and add new P, Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new column = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q: new column = 2
}
